# Added 4wk low sales check - update forecast comparison values and summary metrics

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet updates ---

# Row 2 (W10)
$wsForecast.Range("D2").Value = 42
$wsForecast.Range("H2").Value = 11.52
$wsForecast.Range("L2").Value = 1.19

# Row 3 (W11)
$wsForecast.Range("D3").Value = 43
$wsForecast.Range("H3").Value = 10.32
$wsForecast.Range("L3").Value = 0.98

# Row 4 (W12)
$wsForecast.Range("D4").Value = 46
$wsForecast.Range("H4").Value = 8.84

# Row 5 (W13)
$wsForecast.Range("D5").Value = 48
$wsForecast.Range("H5").Value = 7.39
$wsForecast.Range("L5").Value = 1.15

# Row 6 (W14)
$wsForecast.Range("D6").Value = 50
$wsForecast.Range("H6").Value = 6.19
$wsForecast.Range("L6").Value = 0.9

# Row 7 (W15)
$wsForecast.Range("H7").Value = 5.23
$wsForecast.Range("L7").Value = 0.93

# Row 8 (W16)
$wsForecast.Range("D8").Value = 49
$wsForecast.Range("H8").Value = 4.31
$wsForecast.Range("L8").Value = 1.01

# Row 9 (W17)
$wsForecast.Range("D9").Value = 49
$wsForecast.Range("H9").Value = 3.3
$wsForecast.Range("L9").Value = 1.05

# Row 10 (W18)
$wsForecast.Range("H10").Value = 2.23
$wsForecast.Range("L10").Value = 1.02

# Row 11 (W19)
$wsForecast.Range("H11").Value = 1.19
$wsForecast.Range("J11").Value = "Normal"
$wsForecast.Range("L11").Value = 0.9399999999999999

# Row 12 (W20)
$wsForecast.Range("D12").Value = 52
$wsForecast.Range("H12").Value = 0.19
$wsForecast.Range("L12").Value = 0.99

# Row 13 (W21)
$wsForecast.Range("L13").Value = 0.95

# Row 14 (W22)
$wsForecast.Range("D14").Value = 48
$wsForecast.Range("L14").Value = 0.98

# Row 15 (W23)
$wsForecast.Range("D15").Value = 48
$wsForecast.Range("L15").Value = 0.8

# Row 16 (W24)
$wsForecast.Range("L16").Value = 1.03

# Row 17 (W25)
$wsForecast.Range("L17").Value = 1.1

# --- Summary sheet updates ---
# (Force text format since these columns store numbers as text strings)

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "780"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "381"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "182"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "43"
